$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2700.5532
$ws.Range("I15").Value = 2700.5532
$ws.Range("K15").Value = 8101.659599999999
$ws.Range("M15").Value = -7932.659599999999

# Row 17
$ws.Range("H17").Value = 396.27777
$ws.Range("J17").Value = 388.65714
$ws.Range("L17").Value = 1165.97142
$ws.Range("N17").Value = -1501.97142

# Row 40
$ws.Range("H40").Value = 2333.3333
$ws.Range("I40").Value = 1666.6666
$ws.Range("K40").Value = 1666.6666
$ws.Range("M40").Value = -1491.6666

# Row 62
$ws.Range("H62").Value = 2629.7334
$ws.Range("I62").Value = 1671.4286
$ws.Range("K62").Value = 1671.4286
$ws.Range("M62").Value = -1047.4286

# Row 65
$ws.Range("H65").Value = 2629.7334
$ws.Range("I65").Value = 1671.4286
$ws.Range("K65").Value = 8357.143
$ws.Range("M65").Value = -5237.143

$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 29990
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 29990
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 29990
$ws.Range("M37").Value = ""
$ws.Range("N37").Value = -30536

# Row 110
$ws.Range("H110").Value = 514.9167
$ws.Range("I110").Value = 433.22223
$ws.Range("K110").Value = 433.22223
$ws.Range("M110").Value = 1611.77777

# Row 122
$ws.Range("H122").Value = 2313.6538
$ws.Range("I122").Value = 2355.6667
$ws.Range("J122").Value = 2137.2
$ws.Range("K122").Value = 7067.000100000001
$ws.Range("L122").Value = 6411.599999999999
$ws.Range("M122").Value = -4617.000100000001
$ws.Range("N122").Value = -11311.6

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2870.9473
$ws.Range("I20").Value = 2536.5334
$ws.Range("J20").Value = 4125
$ws.Range("K20").Value = 2536.5334
$ws.Range("L20").Value = 4125
$ws.Range("M20").Value = -2289.5334
$ws.Range("N20").Value = -4619

# Row 99
$ws.Range("H99").Value = 2500.8
$ws.Range("I99").Value = 2399.4
$ws.Range("J99").Value = 2602.2
$ws.Range("K99").Value = 2399.4
$ws.Range("L99").Value = 2602.2
$ws.Range("M99").Value = -901.4000000000001
$ws.Range("N99").Value = -5598.2

# Row 134
$ws.Range("H134").Value = 36389.6
$ws.Range("I134").Value = 49157.363
$ws.Range("J134").Value = 1278.25
$ws.Range("K134").Value = 147472.089
$ws.Range("L134").Value = 3834.75
$ws.Range("M134").Value = -144937.089
$ws.Range("N134").Value = -8904.75

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 63671444
$ws.Range("I6").Value = 9006143
$ws.Range("K6").Value = 9006143
$ws.Range("M6").Value = -9006030

# Row 31
$ws.Range("H31").Value = 19564.646
$ws.Range("I31").Value = 21306.6
$ws.Range("K31").Value = 21306.6
$ws.Range("M31").Value = -21011.6

# Row 34
$ws.Range("H34").Value = 19564.646
$ws.Range("I34").Value = 21306.6
$ws.Range("K34").Value = 21306.6
$ws.Range("M34").Value = -21104.6

# Row 74
$ws.Range("H74").Value = 28080.428
$ws.Range("J74").Value = 30260.5
$ws.Range("L74").Value = 30260.5
$ws.Range("N74").Value = -32008.5

# Row 77
$ws.Range("H77").Value = 28080.428
$ws.Range("J77").Value = 30260.5
$ws.Range("L77").Value = 90781.5
$ws.Range("N77").Value = -99517.5

# Row 86
$ws.Range("H86").Value = 5753778
$ws.Range("I86").Value = 1501.5625
$ws.Range("J86").Value = 12833503
$ws.Range("K86").Value = 1501.5625
$ws.Range("L86").Value = 12833503
$ws.Range("M86").Value = -378.5625
$ws.Range("N86").Value = -12835749

# Row 89
$ws.Range("H89").Value = 5753778
$ws.Range("I89").Value = 1501.5625
$ws.Range("J89").Value = 12833503
$ws.Range("K89").Value = 7507.8125
$ws.Range("L89").Value = 64167515
$ws.Range("M89").Value = -1891.8125
$ws.Range("N89").Value = -64178747

# Row 99
$ws.Range("H99").Value = 21609458
$ws.Range("I99").Value = 5559140.5
$ws.Range("J99").Value = 41672356
$ws.Range("K99").Value = 5559140.5
$ws.Range("L99").Value = 41672356
$ws.Range("M99").Value = -5557642.5
$ws.Range("N99").Value = -41675352

# Row 107
$ws.Range("H107").Value = 1380.3334
$ws.Range("I107").Value = 976.3333
$ws.Range("J107").Value = 1784.3334
$ws.Range("K107").Value = 976.3333
$ws.Range("L107").Value = 1784.3334
$ws.Range("M107").Value = 943.6667
$ws.Range("N107").Value = -5624.3334

# Row 126
$ws.Range("H126").Value = 21609458
$ws.Range("I126").Value = 5559140.5
$ws.Range("J126").Value = 41672356
$ws.Range("K126").Value = 16677421.5
$ws.Range("L126").Value = 125017068
$ws.Range("M126").Value = -16674951.5
$ws.Range("N126").Value = -125022008

# Row 134
$ws.Range("H134").Value = 1176.6538
$ws.Range("I134").Value = 972.4
$ws.Range("J134").Value = 1304.3125
$ws.Range("K134").Value = 2917.2
$ws.Range("L134").Value = 3912.9375
$ws.Range("M134").Value = -382.1999999999998
$ws.Range("N134").Value = -8982.9375

$ws = $wb.Worksheets.Item("CUL")
# Row 21
$ws.Range("H21").Value = 393.33334
$ws.Range("I21").Value = 393.33334
$ws.Range("K21").Value = 1180.00002
$ws.Range("M21").Value = -1007.00002

# Row 24
$ws.Range("H24").Value = 100
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = ""

# Row 25
$ws.Range("H25").Value = 897.25
$ws.Range("I25").Value = 897.25
$ws.Range("K25").Value = 2691.75
$ws.Range("M25").Value = -2522.75

# Row 30
$ws.Range("H30").Value = 897.25
$ws.Range("I30").Value = 897.25
$ws.Range("K30").Value = 2691.75
$ws.Range("M30").Value = -2589.75

# Row 62
$ws.Range("H62").Value = 6730.625
$ws.Range("J62").Value = 7565.2856
$ws.Range("L62").Value = 22695.8568
$ws.Range("N62").Value = -24067.8568

# Row 65
$ws.Range("H65").Value = 6730.625
$ws.Range("J65").Value = 7565.2856
$ws.Range("L65").Value = 68087.5704
$ws.Range("N65").Value = -74951.5704

# Row 69
$ws.Range("H69").Value = 2499.5
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 2499.5
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 7498.5
$ws.Range("M69").Value = ""
$ws.Range("N69").Value = -9120.5

# Row 72
$ws.Range("H72").Value = 2499.5
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 2499.5
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 22495.5
$ws.Range("M72").Value = ""
$ws.Range("N72").Value = -30607.5

# Row 74
$ws.Range("H74").Value = 9325
$ws.Range("I74").Value = 8800
$ws.Range("J74").Value = 9850
$ws.Range("K74").Value = 26400
$ws.Range("L74").Value = 29550
$ws.Range("M74").Value = -25339
$ws.Range("N74").Value = -31672

# Row 77
$ws.Range("H77").Value = 9325
$ws.Range("I77").Value = 8800
$ws.Range("J77").Value = 9850
$ws.Range("K77").Value = 79200
$ws.Range("L77").Value = 88650
$ws.Range("M77").Value = -73896
$ws.Range("N77").Value = -99258

# Row 131
$ws.Range("H131").Value = 775.01
$ws.Range("J131").Value = 785.4316
$ws.Range("L131").Value = 2356.2948
$ws.Range("N131").Value = -12436.2948

# Row 137
$ws.Range("H137").Value = 2355.818
$ws.Range("J137").Value = 4559.8
$ws.Range("L137").Value = 13679.4
$ws.Range("N137").Value = -23879.4

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3857
$ws.Range("I40").Value = 3537.5
$ws.Range("K40").Value = 3537.5
$ws.Range("M40").Value = -3401.5

# Row 118
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = ""

# Row 136
$ws.Range("H136").Value = 126995.25
$ws.Range("J136").Value = 2980
$ws.Range("L136").Value = 8940
$ws.Range("N136").Value = -14040

# Row 140
$ws.Range("H140").Value = 49953.57
$ws.Range("J140").Value = 49953.57
$ws.Range("L140").Value = 49953.57
$ws.Range("N140").Value = -60313.57

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 2080286.8
$ws.Range("I113").Value = 1737.375
$ws.Range("J113").Value = 5405965.5
$ws.Range("K113").Value = 5212.125
$ws.Range("L113").Value = 16217896.5
$ws.Range("M113").Value = -3042.125
$ws.Range("N113").Value = -16222236.5

# Row 126
$ws.Range("H126").Value = 1670.7142
$ws.Range("I126").Value = 1200
$ws.Range("J126").Value = 1706.9231
$ws.Range("K126").Value = 3600
$ws.Range("L126").Value = 5120.7693
$ws.Range("M126").Value = -1130
$ws.Range("N126").Value = -10060.7693
